$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 4449.4346
$ws.Range("I9").Value = 5662.1665
$ws.Range("J9").Value = 83.59999999999999
$ws.Range("K9").Value = 5662.1665
$ws.Range("L9").Value = 83.59999999999999
$ws.Range("M9").Value = -5493.1665
$ws.Range("N9").Value = -421.6
$ws.Range("H12").Value = 14627.429
$ws.Range("I12").Value = 14627.429
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 14627.429
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -14457.429
$ws.Range("N12").Value = ""
$ws.Range("H33").Value = 430.73334
$ws.Range("I33").Value = 457
$ws.Range("K33").Value = 457
$ws.Range("M33").Value = -228
$ws.Range("H58").Value = 201.66667
$ws.Range("I58").Value = 224.375
$ws.Range("J58").Value = 20
$ws.Range("K58").Value = 673.125
$ws.Range("L58").Value = 60
$ws.Range("M58").Value = -523.125
$ws.Range("N58").Value = -360
$ws.Range("H103").Value = 2066
$ws.Range("I103").Value = 1683.2
$ws.Range("K103").Value = 5049.6
$ws.Range("M103").Value = -4463.6
$ws.Range("H113").Value = 4212.25
$ws.Range("I113").Value = 2861.6667
$ws.Range("K113").Value = 2861.6667
$ws.Range("M113").Value = 392.3332999999998
$ws.Range("H127").Value = 1245
$ws.Range("I127").Value = 1245
$ws.Range("K127").Value = 3735
$ws.Range("M127").Value = 1225
$ws.Range("H137").Value = 2202.3262
$ws.Range("I137").Value = 1543.875
$ws.Range("K137").Value = 4631.625
$ws.Range("M137").Value = -2081.625
$ws.Range("H138").Value = 4026.386
$ws.Range("J138").Value = 4010.0212
$ws.Range("L138").Value = 12030.0636
$ws.Range("N138").Value = -22310.0636

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9695.253000000001
$ws.Range("I32").Value = 6500.5083
$ws.Range("J32").Value = 15427
$ws.Range("K32").Value = 6500.5083
$ws.Range("L32").Value = 15427
$ws.Range("M32").Value = -6213.5083
$ws.Range("N32").Value = -16001
$ws.Range("H61").Value = 52635756
$ws.Range("I61").Value = 52635756
$ws.Range("K61").Value = 52635756
$ws.Range("M61").Value = -52635544
$ws.Range("H74").Value = 5142.2334
$ws.Range("I74").Value = 5334.5356
$ws.Range("K74").Value = 5334.5356
$ws.Range("M74").Value = -4460.5356
$ws.Range("H77").Value = 5142.2334
$ws.Range("I77").Value = 5334.5356
$ws.Range("K77").Value = 26672.678
$ws.Range("M77").Value = -22304.678
$ws.Range("H80").Value = 40000
$ws.Range("I80").Value = 30000
$ws.Range("J80").Value = 50000
$ws.Range("K80").Value = 30000
$ws.Range("L80").Value = 50000
$ws.Range("M80").Value = -29002
$ws.Range("N80").Value = -51996
$ws.Range("H83").Value = 40000
$ws.Range("I83").Value = 30000
$ws.Range("J83").Value = 50000
$ws.Range("K83").Value = 90000
$ws.Range("L83").Value = 150000
$ws.Range("M83").Value = -85008
$ws.Range("N83").Value = -159984
$ws.Range("H102").Value = 5292915
$ws.Range("I102").Value = 6174276
$ws.Range("J102").Value = 4750
$ws.Range("K102").Value = 6174276
$ws.Range("L102").Value = 4750
$ws.Range("M102").Value = -6172654
$ws.Range("N102").Value = -7994
$ws.Range("H110").Value = 3988.7
$ws.Range("I110").Value = 3915.2
$ws.Range("J110").Value = 4062.2
$ws.Range("K110").Value = 3915.2
$ws.Range("L110").Value = 4062.2
$ws.Range("M110").Value = -1870.2
$ws.Range("N110").Value = -8152.2
$ws.Range("H122").Value = 2030.9706
$ws.Range("I122").Value = 2147.3872
$ws.Range("J122").Value = 828
$ws.Range("K122").Value = 6442.1616
$ws.Range("L122").Value = 2484
$ws.Range("M122").Value = -3992.1616
$ws.Range("N122").Value = -7384
$ws.Range("H136").Value = 52635756
$ws.Range("I136").Value = 52635756
$ws.Range("K136").Value = 157907268
$ws.Range("M136").Value = -157904718

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2974.8
$ws.Range("I134").Value = 2648.4595
$ws.Range("K134").Value = 7945.3785
$ws.Range("M134").Value = -5410.3785

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1758.7966
$ws.Range("I132").Value = 1745.3928
$ws.Range("J132").Value = 2009
$ws.Range("K132").Value = 5236.178400000001
$ws.Range("L132").Value = 6027
$ws.Range("M132").Value = -2706.178400000001
$ws.Range("N132").Value = -11087
$ws.Range("H134").Value = 1819.1538
$ws.Range("I134").Value = 1877.2727
$ws.Range("K134").Value = 5631.8181
$ws.Range("M134").Value = -3096.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 2593.1538
$ws.Range("J97").Value = 3427.7778
$ws.Range("L97").Value = 10283.3334
$ws.Range("N97").Value = -11275.3334
$ws.Range("H113").Value = 144495.86
$ws.Range("J113").Value = 1898.5
$ws.Range("L113").Value = 5695.5
$ws.Range("N113").Value = -10035.5
$ws.Range("H122").Value = 1794.3684
$ws.Range("J122").Value = 3064.1428
$ws.Range("L122").Value = 27577.2852
$ws.Range("N122").Value = -32477.2852
$ws.Range("H128").Value = 194799.81
$ws.Range("I128").Value = 194799.81
$ws.Range("K128").Value = 584399.4299999999
$ws.Range("M128").Value = -579419.4299999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 60747.19
$ws.Range("I122").Value = 94977.30499999999
$ws.Range("K122").Value = 284931.915
$ws.Range("M122").Value = -282481.915
$ws.Range("H130").Value = 24500
$ws.Range("J130").Value = 24500
$ws.Range("L130").Value = 24500
$ws.Range("N130").Value = -34540

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2053.4614
$ws.Range("I22").Value = 2285.2856
$ws.Range("K22").Value = 2285.2856
$ws.Range("M22").Value = -1990.2856
$ws.Range("H27").Value = 2053.4614
$ws.Range("I27").Value = 2285.2856
$ws.Range("K27").Value = 2285.2856
$ws.Range("M27").Value = -2178.2856
$ws.Range("H55").Value = 588.5
$ws.Range("I55").Value = 441
$ws.Range("K55").Value = 441
$ws.Range("M55").Value = -268
$ws.Range("H68").Value = 1512160.9
$ws.Range("I68").Value = 1636507.6
$ws.Range("K68").Value = 1636507.6
$ws.Range("M68").Value = -1635758.6
$ws.Range("H71").Value = 1512160.9
$ws.Range("I71").Value = 1636507.6
$ws.Range("K71").Value = 8182538
$ws.Range("M71").Value = -8178794
$ws.Range("H100").Value = 10437315
$ws.Range("I100").Value = 11664994
$ws.Range("J100").Value = 2048.5
$ws.Range("K100").Value = 11664994
$ws.Range("L100").Value = 2048.5
$ws.Range("M100").Value = -11664453
$ws.Range("N100").Value = -3130.5
$ws.Range("H122").Value = 4255
$ws.Range("I122").Value = 4297.5
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 12892.5
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -10442.5
$ws.Range("N122").Value = -16900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 16053.375
$ws.Range("I41").Value = 10250
$ws.Range("J41").Value = 17987.834
$ws.Range("K41").Value = 10250
$ws.Range("L41").Value = 17987.834
$ws.Range("M41").Value = -9860
$ws.Range("N41").Value = -18767.834
$ws.Range("H122").Value = 3238.1904
$ws.Range("I122").Value = 2989.0557
$ws.Range("J122").Value = 4733
$ws.Range("K122").Value = 8967.167099999999
$ws.Range("L122").Value = 14199
$ws.Range("M122").Value = -6517.167099999999
$ws.Range("N122").Value = -19099
$ws.Range("H136").Value = 1564.0968
$ws.Range("I136").Value = 1327.138
$ws.Range("K136").Value = 3981.414
$ws.Range("M136").Value = -1431.414
